# "Add walk to the platia and back"
# February's walk distance (column G, the "extra walk" leg) increases by
# 2.8, from 59.2 to 62. Column F ("=F1+G2" running leg total) recalculates
# from 109.2 to 112 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 62

# Bring the data sheet to the front (it was the Chart1 chart-sheet that
# was active/selected before; the edit leaves Sheet1 selected instead).
$ws.Activate()
